# Add a new "AuditNET.AzureStorage.Queue" column to the Versions sheet,
# reflecting the move of AzureQueueDataProvider.cs into its own project and
# its upgrade from Microsoft.Azure.Storage.Queue to Azure.Storage.Queues.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column F (AzureStorage), shifting F:L -> G:M
$ws.Columns("F").Insert()

# Header for the new column
$ws.Range("F1").Value2 = "AuditNET.AzureStorage.Queue"

# Fill the new column with "-" placeholders for rows that have data in this sheet
$ws.Range("F2").Value2  = "-"
$ws.Range("F3").Value2  = "-"
$ws.Range("F4").Value2  = "-"
$ws.Range("F5").Value2  = "-"
$ws.Range("F6").Value2  = "-"
$ws.Range("F7").Value2  = "-"
$ws.Range("F8").Value2  = "-"
$ws.Range("F9").Value2  = "-"
$ws.Range("F10").Value2 = "-"
$ws.Range("F11").Value2 = "-"
$ws.Range("F12").Value2 = "-"
$ws.Range("F13").Value2 = "-"
$ws.Range("F14").Value2 = "-"
$ws.Range("F15").Value2 = "-"
$ws.Range("F16").Value2 = "-"
$ws.Range("F17").Value2 = "-"

# Row 18 (v17.0.0 release) actually shipped the new AzureStorage.Queue package,
# and the AuditNET core version for that release was bumped too.
$ws.Range("E18").Value2 = "6.0.0"
$ws.Range("F18").Value2 = "1.0.0"

# Resize the new column to fit its (longer) header text
$ws.Columns("F").AutoFit()

# Update selection to match author's saved cursor position
$ws.Range("E18").Select()
